# Auto-generated Excel COM-interop script applying the "Horarios actualizados 523" update
# to workbook horarios-141-2026-01-17.xlsx (sheets LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912 - main schedule table
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 15:57:19"
$ws1.Range("A3").Value = "Total filas: 347"

$data1 = @(
    ("05:42:22","07:36","17X38_ROMERO",114,"LP1912"),
    ("06:33:46","07:36","27_EL RETIRO",63,"LP1912"),
    ("05:42:22","07:37","27_EL RETIRO",115,"LP1912"),
    ("07:36:59","07:38","16_SANTA ANA",2,"LP1912"),
    ("06:33:46","07:43","10_OLMOS",70,"LP1912"),
    ("05:57:08","07:44","10_OLMOS",107,"LP1912"),
    ("07:36:59","07:50","15_ABASTO",14,"LP1912"),
    ("05:57:08","07:51","15_ABASTO",114,"LP1912"),
    ("06:16:15","07:58","23_HERNANDEZ",102,"LP1912"),
    ("06:33:46","07:59","11_ETCHEVERRY",86,"LP1912"),
    ("07:12:53","07:59","23_HERNANDEZ",47,"LP1912"),
    ("06:33:46","08:00","23_HERNANDEZ",87,"LP1912"),
    ("06:16:15","08:00","11_ETCHEVERRY",104,"LP1912"),
    ("07:12:53","08:01","16_SANTA ANA",49,"LP1912"),
    ("06:45:50","08:01","23_HERNANDEZ",76,"LP1912"),
    ("06:52:52","08:02","23_HERNANDEZ",70,"LP1912"),
    ("07:36:59","08:02","17_ROMERO",26,"LP1912"),
    ("06:16:15","08:03","17_ROMERO",107,"LP1912"),
    ("06:33:46","08:12","10_OLMOS",99,"LP1912"),
    ("07:12:53","08:13","10_OLMOS",61,"LP1912"),
    ("08:11:27","08:14","11_ETCHEVERRY",3,"LP1912"),
    ("07:36:59","08:14","17_ROMERO",38,"LP1912"),
    ("07:48:35","08:14","10_OLMOS",26,"LP1912"),
    ("06:16:15","08:15","17_ROMERO",119,"LP1912"),
    ("07:36:59","08:25","15X38_ABASTO",49,"LP1912"),
    ("06:33:46","08:26","15X38_ABASTO",113,"LP1912"),
    ("06:33:46","08:27","84_COLONIA URQUIZA-ESC 49",114,"LP1912"),
    ("08:29:19","08:29","23_HERNANDEZ",0,"LP1912"),
    ("06:45:50","08:29","14_ABASTO",104,"LP1912"),
    ("07:36:59","08:30","16_P MOR-SANTA ANA",54,"LP1912"),
    ("06:33:46","08:31","16_P MOR-SANTA ANA",118,"LP1912"),
    ("08:11:27","08:33","23_HERNANDEZ",22,"LP1912"),
    ("06:45:50","08:38","215C_EL PATO",113,"LP1912"),
    ("07:48:35","08:39","215C_EL PATO",51,"LP1912"),
    ("07:12:53","08:43","10_OLMOS",91,"LP1912"),
    ("07:48:35","08:44","10_OLMOS",56,"LP1912"),
    ("07:12:53","08:49","215A_EL PATO",97,"LP1912"),
    ("07:48:35","08:50","215A_EL PATO",62,"LP1912"),
    ("08:11:27","08:53","16_SANTA ANA",42,"LP1912"),
    ("08:29:19","08:54","16_SANTA ANA",25,"LP1912"),
    ("08:53:12","08:55","16_SANTA ANA",2,"LP1912"),
    ("07:12:53","08:59","215B_EL PATO",107,"LP1912"),
    ("08:39:08","09:00","215B_EL PATO",21,"LP1912"),
    ("07:36:59","09:01","17X38_ROMERO",85,"LP1912"),
    ("07:36:59","09:02","23_HERNANDEZ",86,"LP1912"),
    ("07:12:53","09:02","17X38_ROMERO",110,"LP1912"),
    ("08:29:19","09:03","23_HERNANDEZ",34,"LP1912"),
    ("07:36:59","09:04","16_SANTA ANA",88,"LP1912"),
    ("08:39:08","09:04","23_HERNANDEZ",25,"LP1912"),
    ("07:48:35","09:08","16_SANTA ANA",80,"LP1912"),
    ("07:36:59","09:14","15_ABASTO",98,"LP1912"),
    ("07:36:59","09:14","11_ETCHEVERRY",98,"LP1912"),
    ("08:39:08","09:15","11_ETCHEVERRY",36,"LP1912"),
    ("07:36:59","09:16","27_EL RETIRO",100,"LP1912"),
    ("08:39:08","09:17","27_EL RETIRO",38,"LP1912"),
    ("07:55:46","09:21","16_SANTA ANA",86,"LP1912"),
    ("09:21:49","09:22","16_SANTA ANA",1,"LP1912"),
    ("07:36:59","09:26","215_EL PELIGRO",110,"LP1912"),
    ("07:48:35","09:27","215_EL PELIGRO",99,"LP1912"),
    ("07:36:59","09:30","16_P MOR-SANTA ANA",114,"LP1912"),
    ("07:48:35","09:31","16_P MOR-SANTA ANA",103,"LP1912"),
    ("08:39:08","09:32","23_HERNANDEZ",53,"LP1912"),
    ("08:46:25","09:33","23_HERNANDEZ",47,"LP1912"),
    ("07:48:35","09:39","15_ABASTO",111,"LP1912"),
    ("07:48:35","09:44","14_ABASTO",116,"LP1912"),
    ("08:39:08","09:45","14_ABASTO",66,"LP1912"),
    ("07:55:46","09:51","15_ABASTO",116,"LP1912"),
    ("08:53:12","09:54","10_OLMOS",61,"LP1912"),
    ("09:21:49","10:01","16_SANTA ANA",40,"LP1912"),
    ("09:21:49","10:02","215C_EL PATO",41,"LP1912"),
    ("09:21:49","10:03","23_HERNANDEZ",42,"LP1912"),
    ("08:11:27","10:03","215C_EL PATO",112,"LP1912"),
    ("10:04:17","10:04","215C_EL PATO",0,"LP1912"),
    ("08:46:25","10:04","14_ABASTO",78,"LP1912"),
    ("10:04:17","10:05","16_SANTA ANA",1,"LP1912"),
    ("08:39:08","10:05","14_ABASTO",86,"LP1912"),
    ("08:11:27","10:10","10_OLMOS",119,"LP1912"),
    ("08:29:19","10:11","27_EL RETIRO",102,"LP1912"),
    ("08:29:19","10:12","10_OLMOS",103,"LP1912"),
    ("08:29:19","10:14","10_OLMOS",105,"LP1912"),
    ("09:21:49","10:14","17_ROMERO",53,"LP1912"),
    ("08:39:08","10:15","10_OLMOS",96,"LP1912"),
    ("08:29:19","10:15","17_ROMERO",106,"LP1912"),
    ("08:46:25","10:16","10_OLMOS",90,"LP1912"),
    ("08:46:25","10:18","27_EL RETIRO",92,"LP1912"),
    ("09:21:49","10:24","11_ETCHEVERRY",63,"LP1912"),
    ("10:04:17","10:25","16_SANTA ANA",21,"LP1912"),
    ("10:04:17","10:25","11_ETCHEVERRY",21,"LP1912"),
    ("08:29:19","10:26","15X38_ABASTO",117,"LP1912"),
    ("08:39:08","10:30","11_ETCHEVERRY",111,"LP1912"),
    ("10:04:17","10:33","23_HERNANDEZ",29,"LP1912"),
    ("08:39:08","10:34","10_OLMOS",115,"LP1912"),
    ("10:36:18","10:36","10_OLMOS",0,"LP1912"),
    ("08:39:08","10:37","16_P MOR-SANTA ANA",118,"LP1912"),
    ("10:36:18","10:37","16_SANTA ANA",1,"LP1912"),
    ("10:04:17","10:39","15_ABASTO",35,"LP1912"),
    ("10:04:17","10:42","27_EL RETIRO",38,"LP1912"),
    ("09:21:49","10:44","14_ABASTO",83,"LP1912"),
    ("10:04:17","10:45","14_ABASTO",41,"LP1912"),
    ("10:48:14","10:48","10_OLMOS",0,"LP1912"),
    ("10:48:14","10:49","16_SANTA ANA",1,"LP1912"),
    ("10:04:17","10:51","15_ABASTO",47,"LP1912"),
    ("10:36:18","10:54","10_OLMOS",18,"LP1912"),
    ("10:55:25","10:56","16_SANTA ANA",1,"LP1912"),
    ("09:21:49","10:56","27_EL RETIRO",95,"LP1912"),
    ("10:55:25","10:56","10_OLMOS",1,"LP1912"),
    ("10:04:17","10:57","27_EL RETIRO",53,"LP1912"),
    ("10:55:25","10:59","27_EL RETIRO",4,"LP1912"),
    ("09:21:49","11:01","17_ROMERO",100,"LP1912"),
    ("10:36:18","11:03","23_HERNANDEZ",27,"LP1912"),
    ("09:21:49","11:04","14_ABASTO",103,"LP1912"),
    ("10:04:17","11:05","14_ABASTO",61,"LP1912"),
    ("10:36:18","11:11","15_ABASTO",35,"LP1912"),
    ("11:11:31","11:11","16_SANTA ANA",0,"LP1912"),
    ("10:04:17","11:11","23_HERNANDEZ",67,"LP1912"),
    ("09:21:49","11:14","225_C ROCA-H SUR",113,"LP1912"),
    ("09:21:49","11:20","215C_EL PATO",119,"LP1912"),
    ("10:48:14","11:21","10_OLMOS",33,"LP1912"),
    ("10:04:17","11:21","215C_EL PATO",77,"LP1912"),
    ("10:36:18","11:22","10_OLMOS",46,"LP1912"),
    ("10:36:18","11:24","11_ETCHEVERRY",48,"LP1912"),
    ("10:36:18","11:25","16_P MOR-SANTA ANA",49,"LP1912"),
    ("10:04:17","11:25","11_ETCHEVERRY",81,"LP1912"),
    ("10:04:17","11:30","15X38_ABASTO",86,"LP1912"),
    ("10:48:14","11:32","23_HERNANDEZ",44,"LP1912"),
    ("10:36:18","11:33","23_HERNANDEZ",57,"LP1912"),
    ("10:48:14","11:33","16_SANTA ANA",45,"LP1912"),
    ("11:34:25","11:34","16_SANTA ANA",0,"LP1912"),
    ("11:34:25","11:34","23_HERNANDEZ",0,"LP1912"),
    ("10:04:17","11:34","10_OLMOS",90,"LP1912"),
    ("10:36:18","11:35","16_SANTA ANA",59,"LP1912"),
    ("10:04:17","11:37","16_P MOR-SANTA ANA",93,"LP1912"),
    ("10:04:17","11:40","215A_EL PATO",96,"LP1912"),
    ("11:34:25","11:41","10_OLMOS",7,"LP1912"),
    ("10:55:25","11:44","215B_EL PATO",49,"LP1912"),
    ("10:04:17","11:45","215B_EL PATO",101,"LP1912"),
    ("10:55:25","11:53","15_ABASTO",58,"LP1912"),
    ("11:53:59","11:53","225_GOMEZ",0,"LP1912"),
    ("10:04:17","11:54","225_GOMEZ",110,"LP1912"),
    ("11:47:13","11:57","16_SANTA ANA",10,"LP1912"),
    ("11:34:25","12:03","23_HERNANDEZ",29,"LP1912"),
    ("11:53:59","12:04","17_ROMERO",11,"LP1912"),
    ("11:11:31","12:05","17_ROMERO",54,"LP1912"),
    ("11:11:31","12:06","23_HERNANDEZ",55,"LP1912"),
    ("10:48:14","12:07","14_ABASTO",79,"LP1912"),
    ("11:47:13","12:07","23_HERNANDEZ",20,"LP1912"),
    ("12:11:45","12:11","16_SANTA ANA",0,"LP1912"),
    ("12:11:45","12:11","15_ABASTO",0,"LP1912"),
    ("11:53:59","12:16","15_ABASTO",23,"LP1912"),
    ("11:11:31","12:17","15_ABASTO",66,"LP1912"),
    ("11:53:59","12:17","10_OLMOS",24,"LP1912"),
    ("12:11:45","12:17","27_EL RETIRO",6,"LP1912"),
    ("11:11:31","12:18","10_OLMOS",67,"LP1912"),
    ("11:34:25","12:20","17_ROMERO",46,"LP1912"),
    ("11:47:13","12:21","17_ROMERO",34,"LP1912"),
    ("11:53:59","12:28","215C_EL PATO",35,"LP1912"),
    ("10:36:18","12:29","215C_EL PATO",113,"LP1912"),
    ("10:36:18","12:30","11_ETCHEVERRY",114,"LP1912"),
    ("11:53:59","12:30","16_P MOR-SANTA ANA",37,"LP1912"),
    ("10:48:14","12:31","11_ETCHEVERRY",103,"LP1912"),
    ("10:36:18","12:31","16_P MOR-SANTA ANA",115,"LP1912"),
    ("12:32:47","12:32","16_P MOR-SANTA ANA",0,"LP1912"),
    ("12:32:47","12:33","17_ROMERO",1,"LP1912"),
    ("12:32:47","12:33","27_EL RETIRO",1,"LP1912"),
    ("10:55:25","12:36","27_EL RETIRO",101,"LP1912"),
    ("10:48:14","12:37","27_EL RETIRO",109,"LP1912"),
    ("11:53:59","12:39","15X38_ABASTO",46,"LP1912"),
    ("10:48:14","12:40","15X38_ABASTO",112,"LP1912"),
    ("10:55:25","12:42","14_ABASTO",107,"LP1912"),
    ("10:55:25","12:43","15X38_ABASTO",108,"LP1912"),
    ("10:48:14","12:43","14_ABASTO",115,"LP1912"),
    ("12:45:57","12:45","15X38_ABASTO",0,"LP1912"),
    ("12:45:57","12:46","17_ROMERO",1,"LP1912"),
    ("11:53:59","12:50","15_ABASTO",57,"LP1912"),
    ("11:47:13","12:51","15_ABASTO",64,"LP1912"),
    ("12:45:57","12:53","10_OLMOS",8,"LP1912"),
    ("12:32:47","12:54","10_OLMOS",22,"LP1912"),
    ("11:11:31","12:54","15X38_ABASTO",103,"LP1912"),
    ("12:53:14","12:54","17_ROMERO",1,"LP1912"),
    ("11:53:59","13:00","215C_EL PATO",67,"LP1912"),
    ("11:11:31","13:01","215C_EL PATO",110,"LP1912"),
    ("12:32:47","13:03","23_HERNANDEZ",31,"LP1912"),
    ("11:47:13","13:05","23_HERNANDEZ",78,"LP1912"),
    ("11:11:31","13:06","14_ABASTO",115,"LP1912"),
    ("12:11:45","13:06","23_HERNANDEZ",55,"LP1912"),
    ("11:53:59","13:07","23_HERNANDEZ",74,"LP1912"),
    ("11:47:13","13:07","14_ABASTO",80,"LP1912"),
    ("11:53:59","13:10","215_ALUAR",77,"LP1912"),
    ("11:53:59","13:10","16_SANTA ANA",77,"LP1912"),
    ("11:34:25","13:11","215_ALUAR",97,"LP1912"),
    ("11:47:13","13:11","16_SANTA ANA",84,"LP1912"),
    ("13:12:59","13:12","215_ALUAR",0,"LP1912"),
    ("13:12:59","13:13","16_SANTA ANA",1,"LP1912"),
    ("11:34:25","13:18","11_ETCHEVERRY",104,"LP1912"),
    ("11:47:13","13:19","11_ETCHEVERRY",92,"LP1912"),
    ("11:53:59","13:20","16_SANTA ANA",87,"LP1912"),
    ("11:53:59","13:20","17_ROMERO",87,"LP1912"),
    ("12:32:47","13:21","16_SANTA ANA",49,"LP1912"),
    ("11:34:25","13:21","17_ROMERO",107,"LP1912"),
    ("11:53:59","13:29","215A_EL PATO",96,"LP1912"),
    ("11:53:59","13:29","10_OLMOS",96,"LP1912"),
    ("11:47:13","13:30","10_OLMOS",103,"LP1912"),
    ("11:34:25","13:30","215A_EL PATO",116,"LP1912"),
    ("11:53:59","13:30","16_P MOR-SANTA ANA",97,"LP1912"),
    ("11:47:13","13:31","16_P MOR-SANTA ANA",104,"LP1912"),
    ("13:12:59","13:33","23_HERNANDEZ",21,"LP1912"),
    ("12:11:45","13:37","23_HERNANDEZ",86,"LP1912"),
    ("13:39:24","13:39","16_SANTA ANA",0,"LP1912"),
    ("11:53:59","13:39","17X38_ROMERO",106,"LP1912"),
    ("12:32:47","13:39","23_HERNANDEZ",67,"LP1912"),
    ("11:47:13","13:40","17X38_ROMERO",113,"LP1912"),
    ("12:45:57","13:41","23_HERNANDEZ",56,"LP1912"),
    ("12:45:57","13:42","14_ABASTO",57,"LP1912"),
    ("12:53:14","13:42","23_HERNANDEZ",49,"LP1912"),
    ("12:53:14","13:43","14_ABASTO",50,"LP1912"),
    ("13:39:24","13:46","16_SANTA ANA",7,"LP1912"),
    ("12:45:57","13:50","15_ABASTO",65,"LP1912"),
    ("12:32:47","13:51","15_ABASTO",79,"LP1912"),
    ("12:45:57","13:53","225_GOMEZ",68,"LP1912"),
    ("12:11:45","13:54","225_GOMEZ",103,"LP1912"),
    ("12:45:57","14:00","17_ROMERO",75,"LP1912"),
    ("12:11:45","14:01","17_ROMERO",110,"LP1912"),
    ("13:51:48","14:01","16_SANTA ANA",10,"LP1912"),
    ("13:39:24","14:03","23_HERNANDEZ",24,"LP1912"),
    ("13:12:59","14:04","14_ABASTO",52,"LP1912"),
    ("13:12:59","14:10","15_ABASTO",58,"LP1912"),
    ("12:53:14","14:11","15_ABASTO",78,"LP1912"),
    ("13:12:59","14:13","10_OLMOS",61,"LP1912"),
    ("13:39:24","14:14","10_OLMOS",35,"LP1912"),
    ("12:45:57","14:16","27_EL RETIRO",91,"LP1912"),
    ("12:53:14","14:17","27_EL RETIRO",84,"LP1912"),
    ("12:32:47","14:20","215C_EL PATO",108,"LP1912"),
    ("12:53:14","14:21","215C_EL PATO",88,"LP1912"),
    ("12:32:47","14:24","11_ETCHEVERRY",112,"LP1912"),
    ("12:53:14","14:25","11_ETCHEVERRY",92,"LP1912"),
    ("13:51:48","14:25","16_SANTA ANA",34,"LP1912"),
    ("13:12:59","14:31","23_HERNANDEZ",79,"LP1912"),
    ("14:31:57","14:31","15X38_ABASTO",0,"LP1912"),
    ("14:31:57","14:31","16_SANTA ANA",0,"LP1912"),
    ("14:31:57","14:31","15_ABASTO",0,"LP1912"),
    ("13:51:48","14:33","23_HERNANDEZ",42,"LP1912"),
    ("13:39:24","14:34","23_HERNANDEZ",55,"LP1912"),
    ("12:45:57","14:36","16_P MOR-SANTA ANA",111,"LP1912"),
    ("12:53:14","14:37","16_P MOR-SANTA ANA",104,"LP1912"),
    ("13:12:59","14:39","17X38_ROMERO",87,"LP1912"),
    ("12:45:57","14:40","17X38_ROMERO",115,"LP1912"),
    ("14:44:53","14:44","15X38_ABASTO",0,"LP1912"),
    ("12:45:57","14:44","215B_EL PATO",119,"LP1912"),
    ("12:53:14","14:45","215B_EL PATO",112,"LP1912"),
    ("14:52:27","14:53","15X38_ABASTO",1,"LP1912"),
    ("14:52:27","14:54","10_OLMOS",2,"LP1912"),
    ("14:52:27","14:55","16_P MOR-SANTA ANA",3,"LP1912"),
    ("13:12:59","15:00","81_EL PELIGRO",108,"LP1912"),
    ("13:39:24","15:01","81_EL PELIGRO",82,"LP1912"),
    ("14:10:21","15:03","23_HERNANDEZ",53,"LP1912"),
    ("13:51:48","15:04","14_ABASTO",73,"LP1912"),
    ("13:39:24","15:13","14_ABASTO",94,"LP1912"),
    ("14:31:57","15:13","10_OLMOS",42,"LP1912"),
    ("13:39:24","15:14","10_OLMOS",95,"LP1912"),
    ("15:17:21","15:19","11_ETCHEVERRY",2,"LP1912"),
    ("14:52:27","15:20","16_SANTA ANA",28,"LP1912"),
    ("15:17:21","15:23","16_SANTA ANA",6,"LP1912"),
    ("13:39:24","15:24","11_ETCHEVERRY",105,"LP1912"),
    ("13:51:48","15:24","215C_EL PATO",93,"LP1912"),
    ("13:39:24","15:25","16_P MOR-SANTA ANA",106,"LP1912"),
    ("13:39:24","15:25","215C_EL PATO",106,"LP1912"),
    ("15:17:21","15:31","10_OLMOS",14,"LP1912"),
    ("14:10:21","15:33","16_SANTA ANA",83,"LP1912"),
    ("15:17:21","15:33","23_HERNANDEZ",16,"LP1912"),
    ("13:51:48","15:36","17X38_ROMERO",105,"LP1912"),
    ("14:44:53","15:44","15_ABASTO",60,"LP1912"),
    ("13:51:48","15:44","14_ABASTO",113,"LP1912"),
    ("15:45:31","15:45","14_ABASTO",0,"LP1912"),
    ("14:52:27","15:52","15_ABASTO",60,"LP1912"),
    ("14:31:57","15:54","27_EL RETIRO",83,"LP1912"),
    ("15:45:31","15:56","27_EL RETIRO",11,"LP1912"),
    ("14:44:53","15:57","16_SANTA ANA",73,"LP1912"),
    ("14:10:21","15:57","27_EL RETIRO",107,"LP1912"),
    ("15:17:21","15:59","15_ABASTO",42,"LP1912"),
    ("14:31:57","16:00","15_ABASTO",89,"LP1912"),
    ("14:10:21","16:01","15_ABASTO",111,"LP1912"),
    ("15:45:31","16:03","23_HERNANDEZ",18,"LP1912"),
    ("15:57:19","16:04","14_ABASTO",7,"LP1912"),
    ("14:44:53","16:10","10_OLMOS",86,"LP1912"),
    ("14:52:27","16:11","10_OLMOS",79,"LP1912"),
    ("15:17:21","16:14","10_OLMOS",57,"LP1912"),
    ("14:31:57","16:16","10_OLMOS",105,"LP1912"),
    ("14:31:57","16:24","11_ETCHEVERRY",113,"LP1912"),
    ("14:31:57","16:24","215_ALUAR",113,"LP1912"),
    ("14:52:27","16:25","215_ALUAR",93,"LP1912"),
    ("14:44:53","16:31","16_P MOR-SANTA ANA",107,"LP1912"),
    ("14:44:53","16:33","225_GOMEZ",109,"LP1912"),
    ("14:44:53","16:33","23_HERNANDEZ",109,"LP1912"),
    ("14:52:27","16:34","225_GOMEZ",102,"LP1912"),
    ("15:17:21","16:36","23_HERNANDEZ",79,"LP1912"),
    ("14:44:53","16:36","17X38_ROMERO",112,"LP1912"),
    ("14:52:27","16:37","17X38_ROMERO",105,"LP1912"),
    ("15:45:31","16:40","23_HERNANDEZ",55,"LP1912"),
    ("15:57:19","16:41","23_HERNANDEZ",44,"LP1912"),
    ("15:45:31","16:48","16_SANTA ANA",63,"LP1912"),
    ("15:57:19","16:51","16_SANTA ANA",54,"LP1912"),
    ("15:45:31","16:53","10_OLMOS",68,"LP1912"),
    ("15:57:19","16:54","10_OLMOS",57,"LP1912"),
    ("15:17:21","16:57","215C_EL PATO",100,"LP1912"),
    ("15:57:19","17:01","16_SANTA ANA",64,"LP1912"),
    ("15:17:21","17:04","14_ABASTO",107,"LP1912"),
    ("15:45:31","17:07","15_ABASTO",82,"LP1912"),
    ("15:45:31","17:17","17_ROMERO",92,"LP1912"),
    ("15:45:31","17:24","11_ETCHEVERRY",99,"LP1912"),
    ("15:57:19","17:27","15_ABASTO",90,"LP1912"),
    ("15:45:31","17:34","10_OLMOS",109,"LP1912"),
    ("15:45:31","17:35","16_P MOR-SANTA ANA",110,"LP1912"),
    ("15:45:31","17:36","27_EL RETIRO",111,"LP1912"),
    ("15:45:31","17:38","17X38_ROMERO",113,"LP1912"),
    ("15:57:19","17:45","215B_EL PATO",108,"LP1912"),
    ("15:57:19","17:51","215_EL PELIGRO",114,"LP1912")
)

$startRow1 = 37
for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = $startRow1 + $i
    $row = $data1[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215 - secondary schedule table (2 new rows appended)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 15:57:19"
$ws2.Range("A3").Value = "Total filas: 44"

$data2 = @(
    ("15:57:19","17:45","215B_EL PATO",108,"LP1912"),
    ("15:57:19","17:51","215_EL PELIGRO",114,"LP1912")
)

$startRow2 = 48
for ($i = 0; $i -lt $data2.Length; $i++) {
    $r = $startRow2 + $i
    $row = $data2[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173 - only the "last updated" timestamp changes
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 15:57:19"

Write-Output "Applied horarios update to all three sheets."
